$wb = $excel.ActiveWorkbook

# --- 1. Create the new "TimePageTest" sheet by duplicating "AdminPageTest" -------------
# (duplicating keeps the exact bestFit column widths / pageSetup / namespaces that a
#  fresh `Worksheets.Add()` sheet would not have)
$adminSheet = $wb.Worksheets.Item("AdminPageTest")
$adminSheet.Copy($null, $adminSheet)
$timeSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$timeSheet.Name = "TimePageTest"

# Trim the duplicated AdminPageTest data (UserName/Password/User/SearchCriteria/RunMode)
# down to UserName/Password/RunMode with Admin/admin123/Y values.
$timeSheet.Range("C1").Value = "RunMode"
$timeSheet.Range("C2").Value = "Y"
$timeSheet.Range("D1:E2").Clear()

# Leave the cursor where the author's session left it on this sheet.
$timeSheet.Range("E23").Select()

# --- 2. TestSuite sheet: add the TimePageTest row, flip the LoginPageTest RunMode -------
$testSuite = $wb.Worksheets.Item("TestSuite")
$testSuite.Range("B2").Value = "N"
$testSuite.Range("A4").Value = "TimePageTest"
$testSuite.Range("B4").Value = "Y"

# Reuse the bordered style already on row 3 for the freshly added row 4 cells.
$testSuite.Range("A3:B3").Copy()
$testSuite.Range("A4:B4").PasteSpecial(-4122)

# --- 3. LoginPageTest sheet: selection becomes a plain range, no more active tab --------
$loginSheet = $wb.Worksheets.Item("LoginPageTest")
$loginSheet.Range("A1:C2").Select()

# --- 4. Final state: TestSuite is the active tab / cell is A4 --------------------------
$testSuite.Activate()
$testSuite.Range("A4").Select()
